# This workbook's sheet holds weekly Fruta/Hortaliza price records.
# The edit re-shuffles which data row (columns D, L, M, N, O, P, Q, R, S, T)
# lives on which sheet row (rows 2-11 and 13-14; row 12 is untouched).
# We snapshot the "before" values for every affected row first, then
# write them back out according to the row permutation derived from the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary between rows and are involved in the shuffle.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Rows that participate in the permutation (row 12 stays as-is).
$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 13, 14)

# Mapping of source row (in the original workbook) -> destination row
# (where that record's values should end up).
$rowMap = @{
    2  = 7
    3  = 11
    4  = 6
    5  = 14
    6  = 4
    7  = 3
    8  = 9
    9  = 5
    10 = 2
    11 = 13
    13 = 8
    14 = 10
}

# Snapshot all current values before any writes happen, so overlapping
# source/destination rows don't clobber data we still need to read.
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Now write each source row's snapshot into its destination row.
foreach ($r in $rows) {
    $destRow = $rowMap[$r]
    $rowData = $snapshot[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $rowData[$col]
    }
}
